$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Paragraph 1, run 1: date line
Replace-Text 'המאמר היומי של מייק - 26.03.25' 'המאמר היומי של מייק - 24.03.25'

# Paragraph 1, run 2: paper title
Replace-Text 'DoReMi: Optimizing Data Mixtures Speeds Up Language Model Pretraining' 'STAR ATTENTION: EFFICIENT LLM INFERENCE OVER LONG SEQUENCES'

# Paragraph 2
Replace-Text 'המאמר שנסקור היום שייך לתחום שלא הכרתי אז יש סיכוי שיהיו שגיאות בסקירה למרות מאמציי הכבירים. המאמר דן באופטימיזצית אימון של מודלי שפה כאשר יש ברשותנו דאטהסטים מדומיינים שונים. בנימה יותר מתמטית המחברים מציעים שיטה למשקול של הדאטהסטים השונים במהלך האימון. כלומר אם יש לנו d דאטהסטים המטרה היא למצוא וקטור d-מימדי α של מספרים אי שליליים המסתכמים ל-1 כאשר α_k היא ההסתברות לדגום דוגמא מדאטהסט D_k. כלומר אנו מרכיבים את סט האימון שלנו בשני שלבים: בשלב הראשון בוחרים דאטהסט עם דגימה מ- α ולאחר מכן בוחרים דוגמא הדאטהסט הנבחר.' 'הסקירה הזו הולכת להיות קצרה. אפילו מאוד קצרה. המאמר המסוקר מציע שיטה לאופטימיזציה של מנגנון ה-attention בטרנספורמרים עבור מקרה שיש לנו כמה מכונות (נקרא hosts במאמר) להריץ את מודל השפה שלנו. המאמר הוא של חברת אנוידיה דרך אגב וזה דווקא לא מפתיע כי (לפי השמועות 🙂) יש להם די הרבה משאבי חישוב. '

# Paragraph 3
Replace-Text 'דרך אחת פשוטה היא לבחור את α_i פרופורציונלית לגודל הדאטהסט D_i כלומר ככל שדאטהסט גדול מדי הוא ייבחר יותר פעמיים. אפשר לדגום גם בצורה יוניפורמית כאשר כל דאטהסט ייבחר בהסתברות 1/d כאשר d הוא מספר הדאטהסטים. יש שיטות שבוחרות α לפי איכות הדאטהסט ומעדיפים דאטהסטים איכותיים יותר על פני אלו שפחות איכותיים.' 'המודל מחזיר אותי לתקופה העליזה מלפני 4-7 שנים שהייתי עד למבול של מאמרים שהציעו אופטימיזציות שונות למנגנון ה-attention. אתם בטח זוכרים LongFormer, Performer, Reformer, LinFormer וכדומה(שחלקם סקרתי בזמנו) - היה גם Star Transformer דרך אגב. רוב השכלולים שהוצעו בתקופה ההיא דיברו על איך ניתן לזרז את ה-attention בלי לפגוע משמעותית בביצועי המודל - כאשר המודל רץ על מכונה אחת. אז היה מאוד פופולרי האיורים הריבועיים שהיה מצויר בהם הפאטרן של ה-attention כלומר באיזה טוקנים טוקן נתון מתחשב כדי לבנות את ייצוגו ההקשרי (contextualized embedding).'

# Paragraph 4
Replace-Text 'אבל איך לבחור את α בצורה שתמקסם את ביצועי המודל המאומן? זו השאלה שעליה מנסה המאמר לענות. אחת השיטות היא לנסות כל מיני ערכים של α ועבור כל אחד לאמן את המודל (brute-force). עבור מודלים גדולים ומספר גבוה של דאטהסטים d המחיר החישובי (= עלות) עלול להיות עצום. השאלה האם ניתן לעשות משהו חכם מזה?' 'המאמר הזה מציע מנגנון attention שניתן לקרוא לו לוקאלי (מזכיר לי קצת רשתות קונבולוציה על ה-inductive bias שלהם המנצל את התלויות הלוקאלית בתמונות). במאמר זה משהו טיפה יותר מורכב (מזכיר גם LongFormer). כאן מחלקים את חלון הקשר לכמה קבוצות של טוקנים c1,...c_n. כל טוקן בכל קבוצה c_i פרט ל-c1 מחשבת את ה-attention עם הטוקנים בתוך אותה הקבוצה ו-c1 בלבד כאשר טוקנים של c1 מתחשבים בכל הטוקנים לבניית האמבדינג שלהם. כלומר הקבוצה הראשונה של הטוקנים משפיעה על האמבדינגס של כל הטוקנים וגם בעצמה מושפעת מכל הטוקנים בחלון ההקשר. המחברים טוענים שללא הוספה של c1(שזה למעשה התחלת הפרומפט) לכל קבוצות הטוקנים המנגנון סובל מירידה רצינית בביצועים'

# Paragraph 5
Replace-Text 'התשובה על השאלה הזו היא כן וזה מה שהמחברים מציעים. בשלב הראשון המחברים מציעים לאמן מודל M_ref קטן עם α f כלשהו (נגיד יוניפורמי). החברים מציעים להשתמש בשיטת (distributionally robust language modeling (DRO-LM שמאמנת מודל קטן הממזער את השגיאה המקסימלית ביחס מעל כל וקטורי α (החוקיים) יחסית לשגיאה של M_ref (הפרש השגיאות בין M_ref למודל המאומן). השגיאה במקרה הזה היא לוג של הנראות של הטוקן הנכון (עבור כל מודל ממצעים עבור כל הטוקנים עבור כל דאטהסט בנפרד). ' 'כמובן ניתן למקבל את התהליך הזה בקלות בין כמה מכונות (hosts) כאשר כל host מחשב את ה-attention הלוקאלי שלו וגם ה-attention עם c1(בשני שלבים). כל host גם שומר את סכום האקספוננטים של Q ו-K (מכנה של הסופטמקס) עבור הטוקנים שלו. לאחר מכן כל הסכומים האלו מועברים ל-host נוסף שמנרמל את כולם עם סכום אקספוננטים של כל ה-hosts ומחשב את הייצוג הסופי של כל וטוקנים.'

# Paragraph 6
Replace-Text 'אם ראיתם כאן בעיית minimax, אתם צודקים. בצורה איטריבית ממקסמים(כלומר עושים מעלה הגרדיאנט או gradient ascent) את הפרש השגיאות (עבור באצ''ים של דוגמאות) מעל α ולאחר מכן ממזערים את הפרש השגיאות מעל משקלי המודל המאומן (כלומר gradient ascent). וקטור המשקול α הסופי שנבחר על ידי מיצוע של כל וקטורי α עבור כל האיטרציות של בעיית המינימקס הזו. מעניין שהבאצ''ים נדגמים באקראי עבור כל האיטרציות. בשלב האחרון מאמנים מודל גדול עם α שמצאנו בצורה הזו.' 'מנגנון זה מאפשר חישוב ממקובל ומהיר יותר של ה-attention (פחות מכפלות מטריצות) כאשר לטענת המחברים הפגיעה בביצועים לא משמעותית.'

# Remove the "מקווה שהצלחתי..." paragraph entirely (paragraph 7 of 8)
$d.Paragraphs.Item(7).Range.Delete()

# Update the arxiv link (now paragraph 7 of 7)
Replace-Text 'https://arxiv.org/abs/2305.10429' 'https://arxiv.org/abs/2411.17116'

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
